$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2371.7827
$ws.Cells.Item(137, 9).Value = 2333.158
$ws.Cells.Item(137, 11).Value = 6999.474
$ws.Cells.Item(137, 13).Value = -4449.474

$ws.Cells.Item(138, 8).Value = 2924
$ws.Cells.Item(138, 10).Value = 4842.857
$ws.Cells.Item(138, 12).Value = 14528.571
$ws.Cells.Item(138, 14).Value = -24808.571

$ws.Cells.Item(139, 8).Value = 149980.72
$ws.Cells.Item(139, 10).Value = 149980.72
$ws.Cells.Item(139, 12).Value = 149980.72
$ws.Cells.Item(139, 14).Value = -160260.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2663
$ws.Cells.Item(2, 9).Value = 2185.3635
$ws.Cells.Item(2, 11).Value = 2185.3635
$ws.Cells.Item(2, 13).Value = -2072.3635

$ws.Cells.Item(32, 8).Value = 6334.963
$ws.Cells.Item(32, 9).Value = 5232.481
$ws.Cells.Item(32, 11).Value = 5232.481
$ws.Cells.Item(32, 13).Value = -4945.481

$ws.Cells.Item(61, 8).Value = 13006
$ws.Cells.Item(61, 9).Value = 13006
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 13006
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -12794
$ws.Cells.Item(61, 14).ClearContents()

$ws.Cells.Item(80, 8).Value = 42750
$ws.Cells.Item(80, 10).Value = 42750
$ws.Cells.Item(80, 12).Value = 42750
$ws.Cells.Item(80, 14).Value = -44746

$ws.Cells.Item(83, 8).Value = 42750
$ws.Cells.Item(83, 10).Value = 42750
$ws.Cells.Item(83, 12).Value = 128250
$ws.Cells.Item(83, 14).Value = -138234

$ws.Cells.Item(116, 8).Value = 2663
$ws.Cells.Item(116, 9).Value = 2185.3635
$ws.Cells.Item(116, 11).Value = 2185.3635
$ws.Cells.Item(116, 13).Value = 108.6365000000001

$ws.Cells.Item(122, 8).Value = 2097
$ws.Cells.Item(122, 9).Value = 2097
$ws.Cells.Item(122, 11).Value = 6291
$ws.Cells.Item(122, 13).Value = -3841

$ws.Cells.Item(133, 8).Value = 197999.5
$ws.Cells.Item(133, 10).Value = 197999.5
$ws.Cells.Item(133, 12).Value = 197999.5
$ws.Cells.Item(133, 14).Value = -203059.5

$ws.Cells.Item(136, 8).Value = 13006
$ws.Cells.Item(136, 9).Value = 13006
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 39018
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -36468
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2663
$ws.Cells.Item(3, 9).Value = 2185.3635
$ws.Cells.Item(3, 11).Value = 2185.3635
$ws.Cells.Item(3, 13).Value = -2071.3635

$ws.Cells.Item(42, 8).Value = 241633.33
$ws.Cells.Item(42, 10).Value = 241633.33
$ws.Cells.Item(42, 12).Value = 241633.33
$ws.Cells.Item(42, 14).Value = -242289.33

$ws.Cells.Item(95, 8).Value = 93000
$ws.Cells.Item(95, 10).Value = 93000
$ws.Cells.Item(95, 12).Value = 93000
$ws.Cells.Item(95, 14).Value = -98492

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 44242.25
$ws.Cells.Item(28, 10).Value = 44242.25
$ws.Cells.Item(28, 12).Value = 44242.25
$ws.Cells.Item(28, 14).Value = -44732.25

$ws.Cells.Item(31, 8).Value = 83341610
$ws.Cells.Item(31, 9).Value = 142861970
$ws.Cells.Item(31, 11).Value = 142861970
$ws.Cells.Item(31, 13).Value = -142861675

$ws.Cells.Item(34, 8).Value = 83341610
$ws.Cells.Item(34, 9).Value = 142861970
$ws.Cells.Item(34, 11).Value = 142861970
$ws.Cells.Item(34, 13).Value = -142861768

$ws.Cells.Item(132, 8).Value = 3496.4285
$ws.Cells.Item(132, 9).Value = 1928.5
$ws.Cells.Item(132, 11).Value = 5785.5
$ws.Cells.Item(132, 13).Value = -3255.5

$ws.Cells.Item(134, 8).Value = 2750.5938
$ws.Cells.Item(134, 9).Value = 1505.16
$ws.Cells.Item(134, 10).Value = 7198.5713
$ws.Cells.Item(134, 11).Value = 4515.48
$ws.Cells.Item(134, 12).Value = 21595.7139
$ws.Cells.Item(134, 13).Value = -1980.48
$ws.Cells.Item(134, 14).Value = -26665.7139

$ws.Cells.Item(141, 8).Value = 161254.94
$ws.Cells.Item(141, 10).Value = 161254.94
$ws.Cells.Item(141, 12).Value = 161254.94
$ws.Cells.Item(141, 14).Value = -171614.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value = 3751
$ws.Cells.Item(47, 9).Value = 3666.6667
$ws.Cells.Item(47, 11).Value = 11000.0001
$ws.Cells.Item(47, 13).Value = -10569.0001

$ws.Cells.Item(107, 8).Value = 317
$ws.Cells.Item(107, 9).Value = 317
$ws.Cells.Item(107, 11).Value = 951
$ws.Cells.Item(107, 13).Value = 969

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1949.6923
$ws.Cells.Item(122, 9).Value = 1710.6666
$ws.Cells.Item(122, 11).Value = 5131.9998
$ws.Cells.Item(122, 13).Value = -2681.9998

$ws.Cells.Item(132, 8).Value = 7217.05
$ws.Cells.Item(132, 9).Value = 6906.4443
$ws.Cells.Item(132, 11).Value = 20719.3329
$ws.Cells.Item(132, 13).Value = -18189.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 8894.6
$ws.Cells.Item(2, 10).Value = 8888
$ws.Cells.Item(2, 12).Value = 8888
$ws.Cells.Item(2, 14).Value = -9112

$ws.Cells.Item(22, 8).Value = 4066.182
$ws.Cells.Item(22, 9).Value = 2975.8
$ws.Cells.Item(22, 11).Value = 2975.8
$ws.Cells.Item(22, 13).Value = -2680.8

$ws.Cells.Item(27, 8).Value = 4066.182
$ws.Cells.Item(27, 9).Value = 2975.8
$ws.Cells.Item(27, 11).Value = 2975.8
$ws.Cells.Item(27, 13).Value = -2868.8

$ws.Cells.Item(40, 8).Value = 2950.8484
$ws.Cells.Item(40, 9).Value = 1992.5
$ws.Cells.Item(40, 11).Value = 1992.5
$ws.Cells.Item(40, 13).Value = -1856.5

$ws.Cells.Item(46, 8).Value = 7619.96
$ws.Cells.Item(46, 9).Value = 7236.7896
$ws.Cells.Item(46, 11).Value = 7236.7896
$ws.Cells.Item(46, 13).Value = -7048.7896

$ws.Cells.Item(61, 8).Value = 65720.5
$ws.Cells.Item(61, 9).Value = 65720.5
$ws.Cells.Item(61, 11).Value = 65720.5
$ws.Cells.Item(61, 13).Value = -65518.5

$ws.Cells.Item(68, 8).Value = 5333
$ws.Cells.Item(68, 9).Value = 4000
$ws.Cells.Item(68, 11).Value = 4000
$ws.Cells.Item(68, 13).Value = -3251

$ws.Cells.Item(71, 8).Value = 5333
$ws.Cells.Item(71, 9).Value = 4000
$ws.Cells.Item(71, 11).Value = 20000
$ws.Cells.Item(71, 13).Value = -16256

$ws.Cells.Item(82, 8).Value = 1095.381
$ws.Cells.Item(82, 9).Value = 1240.4615
$ws.Cells.Item(82, 11).Value = 1240.4615
$ws.Cells.Item(82, 13).Value = -879.4614999999999

$ws.Cells.Item(85, 8).Value = 1095.381
$ws.Cells.Item(85, 9).Value = 1240.4615
$ws.Cells.Item(85, 11).Value = 1240.4615
$ws.Cells.Item(85, 13).Value = 7.538500000000113

$ws.Cells.Item(113, 8).Value = 65720.5
$ws.Cells.Item(113, 9).Value = 65720.5
$ws.Cells.Item(113, 11).Value = 65720.5
$ws.Cells.Item(113, 13).Value = -63550.5

$ws.Cells.Item(122, 8).Value = 1955.875
$ws.Cells.Item(122, 9).Value = 2092.4285
$ws.Cells.Item(122, 11).Value = 6277.2855
$ws.Cells.Item(122, 13).Value = -3827.2855

$ws.Cells.Item(132, 8).Value = 18176.445
$ws.Cells.Item(132, 9).Value = 22227.285
$ws.Cells.Item(132, 11).Value = 66681.855
$ws.Cells.Item(132, 13).Value = -64151.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 8277.666999999999
$ws.Cells.Item(132, 9).Value = 6987.7896
$ws.Cells.Item(132, 10).Value = 13179.2
$ws.Cells.Item(132, 11).Value = 20963.3688
$ws.Cells.Item(132, 12).Value = 39537.60000000001
$ws.Cells.Item(132, 13).Value = -18433.3688
$ws.Cells.Item(132, 14).Value = -44597.60000000001
